# Update workbook for "Add data for 2021-11-04" commit:
# - Rename sheet + update header label from "October 26" to "October 27"
# - Bump a handful of monthly carjacking counts across neighborhoods,
#   including a few brand-new (previously empty) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name / title update
$ws.Name = "Through 2021-10-27"

# Column header label (B1) text update
$ws.Range("B1").Value = "October 2021 (through October 27)"

# Per-neighborhood / per-month count updates
$ws.Range("V2").Value = 7
$ws.Range("AF2").Value = 8
$ws.Range("BJ2").Value = 4
$ws.Range("B3").Value = 13
$ws.Range("L3").Value = 11
$ws.Range("AP3").Value = 8
$ws.Range("BJ3").Value = 4
$ws.Range("B4").Value = 18
$ws.Range("L4").Value = 11
$ws.Range("L5").Value = 2
$ws.Range("AP6").Value = 3
$ws.Range("L8").Value = 5
$ws.Range("B9").Value = 10
$ws.Range("V9").Value = 5
$ws.Range("V10").Value = 1
$ws.Range("L11").Value = 7
$ws.Range("AP11").Value = 3
$ws.Range("B22").Value = 5
$ws.Range("AP22").Value = 1
$ws.Range("AP24").Value = 2
$ws.Range("B25").Value = 2
$ws.Range("V27").Value = 1
$ws.Range("B38").Value = 6
$ws.Range("AP46").Value = 1
$ws.Range("AP48").Value = 1
$ws.Range("L60").Value = 2
$ws.Range("AP80").Value = 3
$ws.Range("V98").Value = 1

$wb.Save()
